$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values per diff ---

$ws.Range("F4").Value  = "Sprint Aquabike, Aquabike, Super Sprint"
$ws.Range("E5").Value  = "Super Sprint, Sprint, Classic and Ironman 70.3"
$ws.Range("F8").Value  = "Super Sprint, Aquathon, Teams"
$ws.Range("F11").Value = "Sprint Aquabike, Aquabike, Super Sprint"
$ws.Range("E14").Value = "Ironman 70.3, Sprint"
$ws.Range("F14").Value = "Aquabike"
$ws.Range("F15").Value = "Sprint"
$ws.Range("E17").Value = "Long Aqua"
$ws.Range("F17").Value = "Short Aqua"
$ws.Range("E19").Value = "Super Sprint, Sprint, Classic and Ironman 70.3"
$ws.Range("E20").Value = "Sprint, Standard"
$ws.Range("E24").Value = "Ironman 70.3, Sprint"
$ws.Range("F24").Value = "Aquabike"
$ws.Range("E25").Value = "Standard, Aquabike"
$ws.Range("F25").Value = "Sprint"
$ws.Range("E27").Value = "Long Aqua"
$ws.Range("F27").Value = "Short Aqua"
$ws.Range("E28").Value = "Super Sprint, Sprint"
$ws.Range("E29").Value = "Super Sprint, Sprint, Classic and Ironman 70.3"
$ws.Range("E30").Value = "Sprint, Standard"

# --- Append new rows 36-40 ---

$ws.Range("A36").Value = "Hunter League"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "1"
$ws.Range("C36").Value = "Stockton Island"
$ws.Range("D36").Value = "No"
$ws.Range("E36").Value = "Sprint, Standard"
$ws.Range("F36").Value = "N/A"
$ws.Range("G36").Value = "Central Coast Triathlon Club"

$ws.Range("A37").Value = "Hunter League"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "2"
$ws.Range("C37").Value = "Sparke Helmore Triathlon"
$ws.Range("D37").Value = "No"
$ws.Range("E37").Value = "Sprint"
$ws.Range("F37").Value = "Super Sprint"
$ws.Range("G37").Value = "Forster Triathlon Club"

$ws.Range("G38").Value = "Singleton Triathlon Club"
$ws.Range("G39").Value = "Maitland Triathlon Club"
$ws.Range("G40").Value = "Newcastle Traithlon Club"
